# fix bug and run examples notebook
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("withdrawals")

# Fix the amount values (bug fix in source data)
$ws.Range("C3").Value = 2322.07
$ws.Range("C4").Value = 936.41
$ws.Range("C9").Value = 43.62

# Normalize the date number format code to uppercase
$ws.Range("A2:B10").NumberFormat = "YYYY-MM-DD"
